# New Submission Synced: 2026-02-09 17:14:02
# Target sheet: "JSS 3C" (the sheet holding the admissions/AI-score table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3C")

# --- Fix existing row 5: Admission No (C5) was stored as text "1",
#     it should be a genuine number 1 -------------------------------------
$ws.Cells.Item(5, 3).Value = 1

# --- Append the new submission as row 6 -----------------------------------
$ws.Cells.Item(6, 1).Value = "2026-02-09 17:14:02"
$ws.Cells.Item(6, 2).Value = "LAWAN SANI"

# Admission No "18" must stay textual (it looks numeric, so Excel would
# normally coerce it to a number). Force text via NumberFormat, assign it,
# then clear the formatting again so no stray style index is left behind.
$c6 = $ws.Cells.Item(6, 3)
$c6.NumberFormat = "@"
$c6.Value = "18"
$c6.ClearFormats()

$ws.Cells.Item(6, 4).Value = 9
